$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: refreshed "last status check" timestamp
$ws.Range("F1").Value = "Last status check on: 20.02.2022 01:30"

# Row 10 (EuroOil Opustena) price refresh
$ws.Range("B10").Value = 36.5
$ws.Range("C10").Value = 36.9

# D10/E10 switch from numeric cells to plain text cells (no special
# number format, same as an un-styled cell) holding the delta / check
# timestamp as literal text.
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "-0.4"
$ws.Range("E10").Value = "2022-02-20 01:36:47"
$ws.Range("D10:E10").ClearFormats()
